$d = $word.ActiveDocument

# Collect the paragraphs that need to be removed:
#  1) the decorative "────..." separator rule paragraphs
#  2) the empty spacer paragraphs (just <w:pPr><w:spacing w:before="40"/></w:pPr>)
#     that sit directly after a table, before the following section's prose.
$sepText = "────────────────────────────────────────────────────────────"

$targets = New-Object System.Collections.ArrayList

foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "*$sepText*") {
        [void]$targets.Add($p.Range)
    }
    elseif ($t -eq "`r") {
        if ($p.Format.SpaceBefore -eq 2 -and $p.Range.Information(12) -eq $false) {
            [void]$targets.Add($p.Range)
        }
    }
}

# Delete from last to first so earlier ranges stay valid.
for ($i = $targets.Count - 1; $i -ge 0; $i--) {
    $targets[$i].Delete()
}

Write-Output "removed: $($targets.Count)"
